$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 18
$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
$ws.Cells.Item($row, 1).Value = 42625.884247685186
$ws.Cells.Item($row, 2).Value = -30
$ws.Cells.Item($row, 3).Value = 57
$ws.Cells.Item($row, 4).Value = 40
$ws.Cells.Item($row, 5).Value = 50
$ws.Cells.Item($row, 6).Value = 50
$ws.Cells.Item($row, 7).Value = 6914
$ws.Cells.Item($row, 8).Value = 10166
$ws.Cells.Item($row, 9).Value = 1027
$ws.Cells.Item($row, 10).Value = 174
$ws.Cells.Item($row, 11).Value = 124
$ws.Cells.Item($row, 12).Value = 5
$ws.Cells.Item($row, 13).Value = 5
$ws.Cells.Item($row, 14).Value = "Bag"
